$d = $word.ActiveDocument

# Remove the anchored drawing (the "Snip Diagonal Corner Rectangle" shape)
# that lives inside the first paragraph's run.
if ($d.Shapes.Count -gt 0) {
    $d.Shapes.Item(1).Delete()
}

# Insert a new paragraph before the (now shape-less) paragraph and give it
# the new body text.
$firstPara = $d.Paragraphs.Item(1).Range
$firstPara.InsertParagraphBefore()
$d.Paragraphs.Item(1).Range.Text = "On the Insert tab, the galleries include items that are designed to coordinate with the overall look of your document."
